# "added required experimental boolean element to valuesets"
#
# On the "Metadata" sheet:
#   - B7 (the "Experimental" property's value) goes from blank to the
#     literal text "true".
#   - B8 (the "Date" property's value) is refreshed to the new export
#     timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Enter the word with a leading apostrophe so Excel stores it as literal
# text ("true") instead of auto-converting it to the boolean TRUE.
$ws.Range("B7").Value = "'true"

# The apostrophe entry marks the cell with a quote-prefixed style; put the
# row's normal body formatting back by copying it from a sibling cell that
# already carries it, so B7 keeps the same look as the rest of the table.
$ws.Range("B7").ClearFormats()
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)  # xlPasteFormats

# Bump the export Date to the new timestamp.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
